$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "rebill" columns: taxableAmount, vatAmount, cardType, userPaymentId ---

# Headers S1:T1 -> copy header format from B1 (bold, fill, centered) then right-align it
$ws.Range("S1").Value = "taxableAmount"
$ws.Range("T1").Value = "vatAmount"
$ws.Range("B1").Copy()
$ws.Range("S1:T1").PasteSpecial(-4122)
$ws.Range("S1:T1").HorizontalAlignment = -4152

# Headers U1:V1 -> optional/red header style, copied from an existing optional column (I1)
$ws.Range("U1").Value = "cardType"
$ws.Range("V1").Value = "userPaymentId"
$ws.Range("I1").Copy()
$ws.Range("U1:V1").PasteSpecial(-4122)

# Data row 2
$ws.Range("S2").Formula = "=F2/1.16"
$ws.Range("T2").Formula = "=S2*0.16"
$ws.Range("S2:T2").HorizontalAlignment = -4152
$ws.Range("S2:T2").NumberFormat = "0.00"

$ws.Range("U2").Value = "vi"
$ws.Range("V2").Value = "UP_MX_hotgo_95345765_1605768224"
$ws.Range("G2").Copy()
$ws.Range("U2:V2").PasteSpecial(-4122)

# Column widths for the new columns (best effort - engine quantizes to 1/6 increments)
$ws.Columns("S:T").ColumnWidth = 16.33
$ws.Columns("V").ColumnWidth = 36.33

# paymentId value moved out of the missing-record template (N2 cleared)
$ws.Range("N2").ClearContents()

# Clarify the optional-columns note (B8, rich text - keep "rojo" bold & red)
$ws.Range("B8").Value = "las columnas en rojo son opcionales, el resto es obligatorio (no pueden estar vacías, o sea tienen que tener un valor)"
$ws.Range("B8").Characters(17, 4).Font.Bold = $true
$ws.Range("B8").Characters(17, 4).Font.Color = 255

# View: move selection to the newly added columns
$ws.Range("V5").Select()
